$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.023261219476147
$ws.Cells.Item(2, 4).Value = 1.028297207922353
$ws.Cells.Item(2, 5).Value = 0.9926147277508489
$ws.Cells.Item(2, 6).Value = 1.034439382986344
$ws.Cells.Item(2, 9).Value = 1.03255988479546
$ws.Cells.Item(2, 10).Value = 1.028442431729336
$ws.Cells.Item(2, 11).Value = 1.031114434729104
$ws.Cells.Item(2, 12).Value = 0.9955398523336033
$ws.Cells.Item(2, 13).Value = 1.037238853434959
$ws.Cells.Item(2, 14).Value = 1.013551064023137

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.023978425109199
$ws.Cells.Item(3, 4).Value = 1.028811342425153
$ws.Cells.Item(3, 5).Value = 0.9936372048519304
$ws.Cells.Item(3, 6).Value = 1.035434628782692
$ws.Cells.Item(3, 9).Value = 1.032695608688672
$ws.Cells.Item(3, 10).Value = 1.02879924455526
$ws.Cells.Item(3, 11).Value = 1.031437268217526
$ws.Cells.Item(3, 12).Value = 0.9963617723202692
$ws.Cells.Item(3, 13).Value = 1.038042791649186
$ws.Cells.Item(3, 14).Value = 1.013669317856763

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.024443020837871
$ws.Cells.Item(4, 4).Value = 1.02914444909305
$ws.Cells.Item(4, 5).Value = 0.9942998659930995
$ws.Cells.Item(4, 6).Value = 1.036079568322224
$ws.Cells.Item(4, 9).Value = 1.032782505501596
$ws.Cells.Item(4, 10).Value = 1.029029946534061
$ws.Cells.Item(4, 11).Value = 1.031645875306287
$ws.Cells.Item(4, 12).Value = 0.9968940712668345
$ws.Cells.Item(4, 13).Value = 1.038563331816661
$ws.Cells.Item(4, 14).Value = 1.013745763720936

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.024638458531735
$ws.Cells.Item(5, 4).Value = 1.029284587671758
$ws.Cells.Item(5, 5).Value = 0.9945786998346017
$ws.Cells.Item(5, 6).Value = 1.036350926650254
$ws.Cells.Item(5, 9).Value = 1.032818814672741
$ws.Cells.Item(5, 10).Value = 1.029126889478181
$ws.Cells.Item(5, 11).Value = 1.031733503694849
$ws.Cells.Item(5, 12).Value = 0.997117960005301
$ws.Cells.Item(5, 13).Value = 1.038782246582119
$ws.Cells.Item(5, 14).Value = 1.013777883864303

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.024671280461283
$ws.Cells.Item(6, 4).Value = 1.029308123400162
$ws.Cells.Item(6, 5).Value = 0.9946255319796338
$ws.Cells.Item(6, 6).Value = 1.036396502103723
$ws.Cells.Item(6, 9).Value = 1.032824898090585
$ws.Cells.Item(6, 10).Value = 1.029143164002226
$ws.Cells.Item(6, 11).Value = 1.031748212747019
$ws.Cells.Item(6, 12).Value = 0.9971555583673453
$ws.Cells.Item(6, 13).Value = 1.038819007968319
$ws.Cells.Item(6, 14).Value = 1.013783275928384

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.024445631810179
$ws.Cells.Item(7, 4).Value = 1.029146321238772
$ws.Cells.Item(7, 5).Value = 0.9943035907982488
$ws.Cells.Item(7, 6).Value = 1.036083193340251
$ws.Cells.Item(7, 9).Value = 1.032782991540618
$ws.Cells.Item(7, 10).Value = 1.02903124206539
$ws.Cells.Item(7, 11).Value = 1.031647046477711
$ws.Cells.Item(7, 12).Value = 0.9968970624462087
$ws.Cells.Item(7, 13).Value = 1.038566256655049
$ws.Cells.Item(7, 14).Value = 1.013746192981902

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.023503494601282
$ws.Cells.Item(8, 4).Value = 1.028470872564137
$ws.Cells.Item(8, 5).Value = 0.9929600610674301
$ws.Cells.Item(8, 6).Value = 1.034775533609151
$ws.Cells.Item(8, 9).Value = 1.03260594458355
$ws.Cells.Item(8, 10).Value = 1.028563054840753
$ws.Cells.Item(8, 11).Value = 1.031223596716797
$ws.Cells.Item(8, 12).Value = 0.995817528259106
$ws.Cells.Item(8, 13).Value = 1.037510477252163
$ws.Cells.Item(8, 14).Value = 1.01359104318698

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.021847358644587
$ws.Cells.Item(9, 4).Value = 1.027283995419433
$ws.Cells.Item(9, 5).Value = 0.9906006454969559
$ws.Cells.Item(9, 6).Value = 1.032478601613096
$ws.Cells.Item(9, 9).Value = 1.032286904998049
$ws.Cells.Item(9, 10).Value = 1.027736722102389
$ws.Cells.Item(9, 11).Value = 1.030475270560642
$ws.Cells.Item(9, 12).Value = 0.9939188001724441
$ws.Cells.Item(9, 13).Value = 1.035652709598421
$ws.Cells.Item(9, 14).Value = 1.013317113363868

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.020746086084202
$ws.Cells.Item(10, 4).Value = 1.026495100745501
$ws.Cells.Item(10, 5).Value = 0.989033133672735
$ws.Cells.Item(10, 6).Value = 1.03095232806067
$ws.Cells.Item(10, 9).Value = 1.032069504094287
$ws.Cells.Item(10, 10).Value = 1.027185008675492
$ws.Cells.Item(10, 11).Value = 1.02997500652092
$ws.Cells.Item(10, 12).Value = 0.9926553831429383
$ws.Cells.Item(10, 13).Value = 1.034416058847296
$ws.Cells.Item(10, 14).Value = 1.01313415663474

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.020269914219403
$ws.Cells.Item(11, 4).Value = 1.026154081793609
$ws.Cells.Item(11, 5).Value = 0.988355674866747
$ws.Cells.Item(11, 6).Value = 1.030292641862659
$ws.Cells.Item(11, 9).Value = 1.031974258700735
$ws.Cells.Item(11, 10).Value = 1.026945929215685
$ws.Cells.Item(11, 11).Value = 1.029758073740495
$ws.Cells.Item(11, 12).Value = 0.9921088820399291
$ws.Cells.Item(11, 13).Value = 1.033881033305106
$ws.Cells.Item(11, 14).Value = 1.013054859432977

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.020093147765751
$ws.Cells.Item(12, 4).Value = 1.026027500760039
$ws.Cells.Item(12, 5).Value = 0.9881042295826724
$ws.Cells.Item(12, 6).Value = 1.030047786775632
$ws.Cells.Item(12, 9).Value = 1.031938714302162
$ws.Cells.Item(12, 10).Value = 1.026857097973536
$ws.Cells.Item(12, 11).Value = 1.029677449057398
$ws.Cells.Item(12, 12).Value = 0.9919059725120875
$ws.Cells.Item(12, 13).Value = 1.033682370266817
$ws.Cells.Item(12, 14).Value = 1.013025393944772

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.020131059989886
$ws.Cells.Item(13, 4).Value = 1.026054648783681
$ws.Cells.Item(13, 5).Value = 0.9881581567098651
$ws.Cells.Item(13, 6).Value = 1.030100300776279
$ws.Cells.Item(13, 9).Value = 1.031946346198537
$ws.Cells.Item(13, 10).Value = 1.026876153756705
$ws.Cells.Item(13, 11).Value = 1.029694745397204
$ws.Cells.Item(13, 12).Value = 0.9919494934313052
$ws.Cells.Item(13, 13).Value = 1.033724981007322
$ws.Cells.Item(13, 14).Value = 1.013031714882307

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.020255300496967
$ws.Cells.Item(14, 4).Value = 1.026143616748286
$ws.Cells.Item(14, 5).Value = 0.9883348863814464
$ws.Cells.Item(14, 6).Value = 1.030272398349211
$ws.Cells.Item(14, 9).Value = 1.031971323972619
$ws.Cells.Item(14, 10).Value = 1.026938586931189
$ws.Cells.Item(14, 11).Value = 1.0297514102146
$ws.Cells.Item(14, 12).Value = 0.9920921077337197
$ws.Cells.Item(14, 13).Value = 1.033864610335747
$ws.Cells.Item(14, 14).Value = 1.013052424028138

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.020331863139877
$ws.Cells.Item(15, 4).Value = 1.026198444640362
$ws.Cells.Item(15, 5).Value = 0.9884438009545853
$ws.Cells.Item(15, 6).Value = 1.030378457468179
$ws.Cells.Item(15, 9).Value = 1.031986691624027
$ws.Cells.Item(15, 10).Value = 1.026977050595117
$ws.Cells.Item(15, 11).Value = 1.029786317192826
$ws.Cells.Item(15, 12).Value = 0.9921799884222134
$ws.Cells.Item(15, 13).Value = 1.033950649791352
$ws.Cells.Item(15, 14).Value = 1.013065182173824

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.020777702662784
$ws.Cells.Item(16, 4).Value = 1.026517745374706
$ws.Cells.Item(16, 5).Value = 0.9890781214508737
$ws.Cells.Item(16, 6).Value = 1.030996134670188
$ws.Cells.Item(16, 9).Value = 1.032075801899845
$ws.Cells.Item(16, 10).Value = 1.02720087181485
$ws.Cells.Item(16, 11).Value = 1.02998939709416
$ws.Cells.Item(16, 12).Value = 0.9926916645766087
$ws.Cells.Item(16, 13).Value = 1.034451576373009
$ws.Cells.Item(16, 14).Value = 1.013139417767541

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.021057551098079
$ws.Cells.Item(17, 4).Value = 1.026718190217439
$ws.Cells.Item(17, 5).Value = 0.989476357848556
$ws.Cells.Item(17, 6).Value = 1.031383909610422
$ws.Cells.Item(17, 9).Value = 1.032131401871833
$ws.Cells.Item(17, 10).Value = 1.027341220523962
$ws.Cells.Item(17, 11).Value = 1.03011670019451
$ws.Cells.Item(17, 12).Value = 0.9930127773699352
$ws.Cells.Item(17, 13).Value = 1.034765916420545
$ws.Cells.Item(17, 14).Value = 1.013185963802172

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.021220848061268
$ws.Cells.Item(18, 4).Value = 1.026835161916901
$ws.Cells.Item(18, 5).Value = 0.9897087662937556
$ws.Cells.Item(18, 6).Value = 1.031610207891611
$ws.Cells.Item(18, 9).Value = 1.032163725313192
$ws.Cells.Item(18, 10).Value = 1.027423065730645
$ws.Cells.Item(18, 11).Value = 1.030190923454164
$ws.Cells.Item(18, 12).Value = 0.9932001317071769
$ws.Cells.Item(18, 13).Value = 1.034949309213863
$ws.Cells.Item(18, 14).Value = 1.013213105982911

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.021276539262058
$ws.Cells.Item(19, 4).Value = 1.026875055621727
$ws.Cells.Item(19, 5).Value = 0.9897880325774034
$ws.Cells.Item(19, 6).Value = 1.031687389342752
$ws.Cells.Item(19, 9).Value = 1.032174728588681
$ws.Cells.Item(19, 10).Value = 1.027450969773511
$ws.Cells.Item(19, 11).Value = 1.030216226461561
$ws.Cells.Item(19, 12).Value = 0.9932640239640975
$ws.Cells.Item(19, 13).Value = 1.03501184876635
$ws.Cells.Item(19, 14).Value = 1.013222359508985

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.021027519158471
$ws.Cells.Item(20, 4).Value = 1.026696678617377
$ws.Cells.Item(20, 5).Value = 0.9894336180360679
$ws.Cells.Item(20, 6).Value = 1.031342293049667
$ws.Cells.Item(20, 9).Value = 1.032125447592291
$ws.Cells.Item(20, 10).Value = 1.02732616427076
$ws.Cells.Item(20, 11).Value = 1.030103044921476
$ws.Cells.Item(20, 12).Value = 0.9929783193494215
$ws.Cells.Item(20, 13).Value = 1.034732186193108
$ws.Cells.Item(20, 14).Value = 1.013180970608836

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.020218711851457
$ws.Cells.Item(21, 4).Value = 1.026117415440226
$ws.Cells.Item(21, 5).Value = 0.9882828385668249
$ws.Cells.Item(21, 6).Value = 1.030221714875196
$ws.Cells.Item(21, 9).Value = 1.031963973218397
$ws.Cells.Item(21, 10).Value = 1.026920202639365
$ws.Cells.Item(21, 11).Value = 1.029734725103974
$ws.Cells.Item(21, 12).Value = 0.9920501090198102
$ws.Cells.Item(21, 13).Value = 1.033823491056015
$ws.Cells.Item(21, 14).Value = 1.013046326000649

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.019710790885551
$ws.Cells.Item(22, 4).Value = 1.025753723127697
$ws.Cells.Item(22, 5).Value = 0.9875604150241495
$ws.Cells.Item(22, 6).Value = 1.029518215088904
$ws.Cells.Item(22, 9).Value = 1.031861487661395
$ws.Cells.Item(22, 10).Value = 1.026664805776049
$ws.Cells.Item(22, 11).Value = 1.029502881149922
$ws.Cells.Item(22, 12).Value = 0.9914670000341481
$ws.Cells.Item(22, 13).Value = 1.033252559780967
$ws.Cells.Item(22, 14).Value = 1.012961606281245

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.019979991122759
$ws.Cells.Item(23, 4).Value = 1.025946474034931
$ws.Cells.Item(23, 5).Value = 0.9879432794643023
$ws.Cells.Item(23, 6).Value = 1.029891053430611
$ws.Cells.Item(23, 9).Value = 1.031915907993207
$ws.Cells.Item(23, 10).Value = 1.026800210575549
$ws.Cells.Item(23, 11).Value = 1.029625810907178
$ws.Cells.Item(23, 12).Value = 0.991776070289318
$ws.Cells.Item(23, 13).Value = 1.033555182765663
$ws.Cells.Item(23, 14).Value = 1.013006523674832

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.021041089095257
$ws.Cells.Item(24, 4).Value = 1.026706398611894
$ws.Cells.Item(24, 5).Value = 0.9894529299347244
$ws.Cells.Item(24, 6).Value = 1.031361097424693
$ws.Cells.Item(24, 9).Value = 1.032128138405864
$ws.Cells.Item(24, 10).Value = 1.027332967598705
$ws.Cells.Item(24, 11).Value = 1.03010921524592
$ws.Cells.Item(24, 12).Value = 0.9929938892766442
$ws.Cells.Item(24, 13).Value = 1.034747427296518
$ws.Cells.Item(24, 14).Value = 1.013183226840687

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.022275020718611
$ws.Cells.Item(25, 4).Value = 1.027590423855988
$ws.Cells.Item(25, 5).Value = 0.9912096547607049
$ws.Cells.Item(25, 6).Value = 1.033071536013649
$ws.Cells.Item(25, 9).Value = 1.032370217031989
$ws.Cells.Item(25, 10).Value = 1.027950499245181
$ws.Cells.Item(25, 11).Value = 1.030668979006
$ws.Cells.Item(25, 12).Value = 0.9944092447426414
$ws.Cells.Item(25, 13).Value = 1.036132664845238
$ws.Cells.Item(25, 14).Value = 1.013387991884234
